$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2025-06-12 Thursday" "2025-06-13 Friday"

Replace-Text "47×17=799" "30×89=2670"
Replace-Text "50×26=1300" "84×13=1092"
Replace-Text "87×31=2697" "36×53=1908"
Replace-Text "39×90=3510" "89×90=8010"
Replace-Text "17×35=595" "84×61=5124"

Replace-Text "54×57=3078" "56×47=2632"
Replace-Text "85×66=5610" "95×59=5605"
Replace-Text "59×71=4189" "93×70=6510"
Replace-Text "73×50=3650" "20×99=1980"
Replace-Text "45×24=1080" "78×28=2184"

Replace-Text "90×69=6210" "94×86=8084"
Replace-Text "29×44=1276" "15×68=1020"
Replace-Text "95×97=9215" "50×90=4500"
Replace-Text "42×53=2226" "28×91=2548"
Replace-Text "47×94=4418" "32×68=2176"

Replace-Text "19×40=760" "43×27=1161"
Replace-Text "65×45=2925" "39×25=975"
Replace-Text "58×23=1334" "89×29=2581"
Replace-Text "45×47=2115" "55×27=1485"
Replace-Text "50×44=2200" "91×91=8281"

Replace-Text "67×73=4891" "23×63=1449"
Replace-Text "86×28=2408" "36×60=2160"
Replace-Text "28×55=1540" "14×45=630"
Replace-Text "83×51=4233" "22×53=1166"
Replace-Text "65×60=3900" "50×33=1650"
